# Plot addition for assignment completion
# Restructure the "Count" (e.g. "10 out of 40") text column into two
# numeric columns (completed / total) so the data can be charted.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row - lower-cased per the target layout, split "Count" into
# "completed" / "total".
$ws.Range("A1").Value = "session"
$ws.Range("B1").Value = "completed"
$ws.Range("C1").Value = "total"

# Data rows: split the old "X out of Y" strings into numeric completed/total.
$ws.Range("B2").Value = 10
$ws.Range("C2").Value = 40

$ws.Range("B3").Value = 10
$ws.Range("C3").Value = 38

$ws.Range("B4").Value = 15
$ws.Range("C4").Value = 42

$ws.Range("B5").Value = 22
$ws.Range("C5").Value = 48

# Select the header row, mirroring the resulting selection in the workbook
# (selecting the full row 1 before inserting a chart from the data).
$ws.Range("A1:XFD1").Select()
